{"js": "// The document contains one table whose five \"data\" rows (0, 4, 8, 12, 16)\n// each hold five two-digit-division fact cells; the other rows are blank\n// work rows. The edit swaps the 25 division facts for new ones while\n// leaving everything else (formatting, layout, row/column count) intact.\n// We address each target cell positionally (row index, column index) and\n// replace its text in place so the existing run/paragraph formatting\n// (font, size, justification) is preserved.\n\nconst newValues = [\n  [0, 0, \"44\u00f75=8, 4\"],\n  [0, 1, \"28\u00f74=7, 0\"],\n  [0, 2, \"20\u00f73=6, 2\"],\n  [0, 3, \"75\u00f75=15, 0\"],\n  [0, 4, \"21\u00f74=5, 1\"],\n  [4, 0, \"22\u00f77=3, 1\"],\n  [4, 1, \"27\u00f74=6, 3\"],\n  [4, 2, \"84\u00f74=21, 0\"],\n  [4, 3, \"31\u00f75=6, 1\"],\n  [4, 4, \"79\u00f74=19, 3\"],\n  [8, 0, \"64\u00f72=32, 0\"],\n  [8, 1, \"14\u00f73=4, 2\"],\n  [8, 2, \"47\u00f72=23, 1\"],\n  [8, 3, \"36\u00f72=18, 0\"],\n  [8, 4, \"59\u00f73=19, 2\"],\n  [12, 0, \"36\u00f77=5, 1\"],\n  [12, 1, \"40\u00f75=8, 0\"],\n  [12, 2, \"14\u00f76=2, 2\"],\n  [12, 3, \"21\u00f78=2, 5\"],\n  [12, 4, \"92\u00f72=46, 0\"],\n  [16, 0, \"92\u00f79=10, 2\"],\n  [16, 1, \"50\u00f77=7, 1\"],\n  [16, 2, \"88\u00f73=29, 1\"],\n  [16, 3, \"37\u00f72=18, 1\"],\n  [16, 4, \"24\u00f79=2, 6\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [row, col, text] of newValues) {\n  const cell = table.getCell(row, col);\n  cell.body.getRange().insertText(text, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains one table whose five \"data\" rows (1, 5, 9, 13, 17 in\n# Word's 1-based indexing) each hold five two-digit-division fact cells; the\n# other rows are blank work rows. The edit swaps the 25 division facts for\n# new ones while leaving everything else (formatting, layout, row/column\n# count) intact. Each target cell is addressed positionally (row, column)\n# and its Range.Text is replaced in place so the existing paragraph/run\n# formatting (font, size, justification) is preserved.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n  \"1,1\"  = \"44\u00f75=8, 4\"\n  \"1,2\"  = \"28\u00f74=7, 0\"\n  \"1,3\"  = \"20\u00f73=6, 2\"\n  \"1,4\"  = \"75\u00f75=15, 0\"\n  \"1,5\"  = \"21\u00f74=5, 1\"\n  \"5,1\"  = \"22\u00f77=3, 1\"\n  \"5,2\"  = \"27\u00f74=6, 3\"\n  \"5,3\"  = \"84\u00f74=21, 0\"\n  \"5,4\"  = \"31\u00f75=6, 1\"\n  \"5,5\"  = \"79\u00f74=19, 3\"\n  \"9,1\"  = \"64\u00f72=32, 0\"\n  \"9,2\"  = \"14\u00f73=4, 2\"\n  \"9,3\"  = \"47\u00f72=23, 1\"\n  \"9,4\"  = \"36\u00f72=18, 0\"\n  \"9,5\"  = \"59\u00f73=19, 2\"\n  \"13,1\" = \"36\u00f77=5, 1\"\n  \"13,2\" = \"40\u00f75=8, 0\"\n  \"13,3\" = \"14\u00f76=2, 2\"\n  \"13,4\" = \"21\u00f78=2, 5\"\n  \"13,5\" = \"92\u00f72=46, 0\"\n  \"17,1\" = \"92\u00f79=10, 2\"\n  \"17,2\" = \"50\u00f77=7, 1\"\n  \"17,3\" = \"88\u00f73=29, 1\"\n  \"17,4\" = \"37\u00f72=18, 1\"\n  \"17,5\" = \"24\u00f79=2, 6\"\n}\n\n$rows = @(1, 5, 9, 13, 17)\nforeach ($r in $rows) {\n  for ($c = 1; $c -le 5; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[\"$r,$c\"]\n  }\n}\n"}
